$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Add a "Code" column to the header row of the area inventory
# report table (row 9), and drop the now-unused spacer/leftover
# rows 10-11 that used to hold that text off to the side.
# ---------------------------------------------------------------

# 1) Break the old A9:B9 merge ("Product Description" used to span
#    both cells) so A9 and B9 become independent header cells.
$ws.Range("A9:B9").UnMerge()

# 2) Slide the existing "Product Description" header into B9 and
#    put the new "Code" header in A9.
$ws.Range("B9").Value = $ws.Range("A9").Text
$ws.Range("A9").Value = "Code"

# 3) Give the new/moved header cells the same bordered, centered
#    look as the rest of the header row (copy format from C9).
$ws.Range("C9").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# 4) Remove the old blank row 10 and the stray row 11 (which used
#    to hold the "Code" text off in column E) entirely - everything
#    needed now lives in the row 9 header.
$ws.Range("A10:E11").EntireRow.Delete()

# ---------------------------------------------------------------
# Widen column A slightly and column B a lot now that it carries
# the full "Product Description" header/values.
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 14.333333333333334
$ws.Columns.Item(2).ColumnWidth = 40.0

# ---------------------------------------------------------------
# Match the updated selection from the source workbook.
# ---------------------------------------------------------------
$ws.Range("I11").Select() | Out-Null
